$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.093.75'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.18%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.844.17'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.06%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.29'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4525'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.79%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3894'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.63'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.64%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07757'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9712'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.25'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.40%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.842.22'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.755'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.27%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.927'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.85%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.86'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.10%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06530'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001013'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.85'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.017'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.62%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.098.35'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.247'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.55'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.98%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.242'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.064.42'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.08'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.03'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.68%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.220'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.98%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '115.90'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.86%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09214'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.93%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.599'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.59%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.365'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.147'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.20%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05978'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.99%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02175'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.81%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.070'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.163'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.004'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5619'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.47%  '

$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.903'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.17%  '

$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1775'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.238'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.20%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.260'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +22.38%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.43%  '

$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5331'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.25%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.67'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.94%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.857'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -5.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '108.92'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.42%  '
